$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173456311225891
$ws.Range("B1").Value = 2.390925168991089
$ws.Range("D1").Value = 2.36376953125
$ws.Range("E1").Value = 1.209546804428101
